# Add data for 2025-05-11
#
# The source diff updates the "2025" column (column L) across the
# "Citywide Totals" sheet, the "By Neighborhood" roll-up sheet, and every
# individual neighborhood sheet that recorded violent-crime incidents on
# 2025-05-11. Each L-column cell receives its new year-to-date total
# (previous total + the day's new incidents for that crime
# category/neighborhood); "Total" row/column cells are recomputed sums.
#
# $wb (ActiveWorkbook) is already open; update each affected worksheet by
# name and write the new cumulative values into column L.

$wb = $excel.ActiveWorkbook

## SHEET: Citywide Totals (rId1)
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 2189
$ws.Range('L3').Value = 2207
$ws.Range('L4').Value = 611
$ws.Range('L5').Value = 130
$ws.Range('L6').Value = 1994
$ws.Range('L7').Value = 7131

## SHEET: By Neighborhood (rId2)
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L2').Value = 56
$ws.Range('L5').Value = 23
$ws.Range('L7').Value = 232
$ws.Range('L8').Value = 451
$ws.Range('L9').Value = 44
$ws.Range('L10').Value = 47
$ws.Range('L15').Value = 50
$ws.Range('L18').Value = 51
$ws.Range('L19').Value = 203
$ws.Range('L20').Value = 184
$ws.Range('L25').Value = 38
$ws.Range('L27').Value = 73
$ws.Range('L29').Value = 365
$ws.Range('L32').Value = 12
$ws.Range('L33').Value = 318
$ws.Range('L34').Value = 46
$ws.Range('L37').Value = 255
$ws.Range('L42').Value = 221
$ws.Range('L47').Value = 52
$ws.Range('L48').Value = 98
$ws.Range('L49').Value = 41
$ws.Range('L52').Value = 142
$ws.Range('L53').Value = 90
$ws.Range('L54').Value = 147
$ws.Range('L55').Value = 64
$ws.Range('L60').Value = 41
$ws.Range('L65').Value = 137
$ws.Range('L67').Value = 254
$ws.Range('L68').Value = 20
$ws.Range('L69').Value = 19
$ws.Range('L72').Value = 31
$ws.Range('L75').Value = 29
$ws.Range('L76').Value = 76
$ws.Range('L79').Value = 197
$ws.Range('L83').Value = 173
$ws.Range('L84').Value = 72
$ws.Range('L85').Value = 376
$ws.Range('L88').Value = 101
$ws.Range('L93').Value = 38
$ws.Range('L94').Value = 83
$ws.Range('L95').Value = 104
$ws.Range('L96').Value = 66
$ws.Range('L99').Value = 111
$ws.Range('L101').Value = 7131

## SHEET: West Ridge (rId4)
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L3').Value = 13
$ws.Range('L6').Value = 16
$ws.Range('L7').Value = 66

## SHEET: Auburn Gresham (rId5)
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L6').Value = 65
$ws.Range('L7').Value = 232

## SHEET: South Shore (rId8)
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L2').Value = 114
$ws.Range('L3').Value = 156
$ws.Range('L6').Value = 69
$ws.Range('L7').Value = 376

## SHEET: Little Village (rId9)
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('L3').Value = 43
$ws.Range('L5').Value = 3
$ws.Range('L7').Value = 142

## SHEET: Norwood Park (rId10)
$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range('L2').Value = 9
$ws.Range('L7').Value = 19

## SHEET: Logan Square (rId11)
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('L2').Value = 31
$ws.Range('L7').Value = 90

## SHEET: Austin (rId12)
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 127
$ws.Range('L3').Value = 152
$ws.Range('L6').Value = 120
$ws.Range('L7').Value = 451

## SHEET: South Chicago (rId13)
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L2').Value = 54
$ws.Range('L5').Value = 3
$ws.Range('L6').Value = 41
$ws.Range('L7').Value = 173

## SHEET: Garfield Park (rId14)
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 82
$ws.Range('L6').Value = 110
$ws.Range('L7').Value = 318

## SHEET: West Pullman (rId15)
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('L2').Value = 42
$ws.Range('L7').Value = 104

## SHEET: Grand Crossing (rId16)
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L3').Value = 75
$ws.Range('L4').Value = 18
$ws.Range('L5').Value = 9
$ws.Range('L6').Value = 74
$ws.Range('L7').Value = 255

## SHEET: New City (rId17)
$ws = $wb.Worksheets.Item('New City')
$ws.Range('L2').Value = 51
$ws.Range('L7').Value = 137

## SHEET: Woodlawn (rId18)
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('L3').Value = 50
$ws.Range('L7').Value = 111

## SHEET: North Lawndale (rId21)
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L2').Value = 74
$ws.Range('L3').Value = 83
$ws.Range('L4').Value = 23
$ws.Range('L6').Value = 68
$ws.Range('L7').Value = 254

## SHEET: South Deering (rId22)
$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('L3').Value = 29
$ws.Range('L7').Value = 72

## SHEET: Lincoln Park (rId23)
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('L6').Value = 19
$ws.Range('L7').Value = 41

## SHEET: Loop (rId24)
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L6').Value = 74
$ws.Range('L7').Value = 147

## SHEET: Englewood (rId25)
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L3').Value = 128
$ws.Range('L7').Value = 365

## SHEET: Lake View (rId26)
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('L6').Value = 36
$ws.Range('L7').Value = 98

## SHEET: Chatham (rId27)
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L2').Value = 66
$ws.Range('L3').Value = 65
$ws.Range('L7').Value = 203

## SHEET: River North (rId29)
$ws = $wb.Worksheets.Item('River North')
$ws.Range('L2').Value = 14
$ws.Range('L7').Value = 76

## SHEET: Humboldt Park (rId32)
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L2').Value = 56
$ws.Range('L7').Value = 221

## SHEET: Avondale (rId34)
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('L2').Value = 22
$ws.Range('L7').Value = 47

## SHEET: Lower West Side (rId36)
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('L3').Value = 23
$ws.Range('L6').Value = 12
$ws.Range('L7').Value = 64

## SHEET: Roseland (rId42)
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L3').Value = 71
$ws.Range('L6').Value = 41
$ws.Range('L7').Value = 197

## SHEET: Chicago Lawn (rId44)
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L3').Value = 59
$ws.Range('L7').Value = 184

## SHEET: Calumet Heights (rId45)
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('L6').Value = 8
$ws.Range('L7').Value = 51

## SHEET: West Lawn (rId48)
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('L4').Value = 2
$ws.Range('L7').Value = 38

## SHEET: Garfield Ridge (rId50)
$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('L4').Value = 5
$ws.Range('L7').Value = 46

## SHEET: West Loop (rId51)
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('L4').Value = 12
$ws.Range('L7').Value = 83

## SHEET: East Side (rId52)
$ws = $wb.Worksheets.Item('East Side')
$ws.Range('L2').Value = 12
$ws.Range('L3').Value = 20
$ws.Range('L7').Value = 38

## SHEET: Kenwood (rId53)
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('L3').Value = 18
$ws.Range('L7').Value = 52

## SHEET: Brighton Park (rId54)
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('L6').Value = 11
$ws.Range('L7').Value = 50

## SHEET: Avalon Park (rId61)
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('L3').Value = 20
$ws.Range('L7').Value = 44

## SHEET: Albany Park (rId64)
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('L2').Value = 15
$ws.Range('L3').Value = 17
$ws.Range('L6').Value = 19
$ws.Range('L7').Value = 56

## SHEET: United Center (rId68)
$ws = $wb.Worksheets.Item('United Center')
$ws.Range('L2').Value = 27
$ws.Range('L7').Value = 101

## SHEET: Galewood (rId69)
$ws = $wb.Worksheets.Item('Galewood')
$ws.Range('L2').Value = 7
$ws.Range('L7').Value = 12

## SHEET: Armour Square (rId70)
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('L4').Value = 1
$ws.Range('L7').Value = 23

## SHEET: Edgewater (rId71)
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('L5').Value = 2
$ws.Range('L7').Value = 73

## SHEET: Pullman (rId73)
$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('L6').Value = 2
$ws.Range('L7').Value = 29

## SHEET: North Park (rId76)
$ws = $wb.Worksheets.Item('North Park')
$ws.Range('L6').Value = 5
$ws.Range('L7').Value = 20

## SHEET: Morgan Park (rId78)
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('L6').Value = 11
$ws.Range('L7').Value = 41

## SHEET: Old Town (rId82)
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('L3').Value = 6
$ws.Range('L7').Value = 31

